# The upstream NATMI TPM-recompute changed the average/total expression
# value of Cxadr in the "ECs" cluster (columns G/H act as both the ligand-
# and receptor-side average/total expression value, depending on row).
# Every other touched cell (I/J/M/N/O/P specificity + Q/R/S/T edge weights)
# is a pure function of the three per-cluster average expression values, so
# recompute them all from that single new number instead of hand-copying
# every cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-cluster ligand/receptor average expression value (Cxadr-Cxadr sheet).
# Only "ECs" changed with the new TPM numbers; FAPs/MuSCs are unchanged.
$avg = @{
    "ECs"   = 3.584831333333333
    "FAPs"  = 0.8953543333333333
    "MuSCs" = 0.8181726666666668
}
$nCells = 3

$total = @{}
foreach ($k in $avg.Keys) { $total[$k] = $avg[$k] * $nCells }

$sumAvg = 0.0
$sumTotal = 0.0
foreach ($k in $avg.Keys) { $sumAvg += $avg[$k]; $sumTotal += $total[$k] }

$specAvg = @{}
$specTotal = @{}
foreach ($k in $avg.Keys) {
    $specAvg[$k]   = $avg[$k]   / $sumAvg
    $specTotal[$k] = $total[$k] / $sumTotal
}

# row -> (sending cluster, target cluster), matching sheet order
$rows = @(
    @{ Row = 2;  Send = "ECs";   Target = "ECs"   },
    @{ Row = 3;  Send = "ECs";   Target = "FAPs"  },
    @{ Row = 4;  Send = "ECs";   Target = "MuSCs" },
    @{ Row = 5;  Send = "FAPs";  Target = "ECs"   },
    @{ Row = 6;  Send = "FAPs";  Target = "FAPs"  },
    @{ Row = 7;  Send = "FAPs";  Target = "MuSCs" },
    @{ Row = 8;  Send = "MuSCs"; Target = "ECs"   },
    @{ Row = 9;  Send = "MuSCs"; Target = "FAPs"  },
    @{ Row = 10; Send = "MuSCs"; Target = "MuSCs" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $send = $r.Send
    $target = $r.Target

    $G = $avg[$send]
    $H = $total[$send]
    $I = $specAvg[$send]
    $J = $specTotal[$send]

    $M = $avg[$target]
    $N = $total[$target]
    $O = $specAvg[$target]
    $P = $specTotal[$target]

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Range("G$row").Value = $G
    $ws.Range("H$row").Value = $H
    $ws.Range("I$row").Value = $I
    $ws.Range("J$row").Value = $J
    $ws.Range("M$row").Value = $M
    $ws.Range("N$row").Value = $N
    $ws.Range("O$row").Value = $O
    $ws.Range("P$row").Value = $P
    $ws.Range("Q$row").Value = $Q
    $ws.Range("R$row").Value = $R
    $ws.Range("S$row").Value = $S
    $ws.Range("T$row").Value = $T
}
